$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill "Usuario" (D) down for the new regression rows
$ws.Range("D2").Value = "su"
$ws.Range("D3").Value = "su"
$ws.Range("D4").Value = "su"

# NroCuenta (F) - unique account numbers per row
$ws.Range("F2").Value = 2240451788
$ws.Range("F3").Value = 7821451462
$ws.Range("F4").Value = 6759658789

# Fill "FechaInicio" (J) down for the new regression rows
$ws.Range("J2").Value = "'07/04/2021"
$ws.Range("J3").Value = "'07/04/2021"
$ws.Range("J4").Value = "'07/04/2021"

# Row 3 payment method details (Tarjeta de Credito)
$ws.Range("K3").Value = "Tarjeta de Crédito"
$ws.Range("L3").Value = "Sí"
$ws.Range("M3").Value = "6 Cuotas - ARS"

# Row 4 payment method (Debito Bancario)
$ws.Range("K4").Value = "Débito Bancario"

# Fill "Producto" (X) down for the new regression rows
$ws.Range("X2").Value = "TR - Todo Riesgo Franquicia Fija"
$ws.Range("X3").Value = "TR - Todo Riesgo Franquicia Fija"
$ws.Range("X4").Value = "TR - Todo Riesgo Franquicia Fija"

# Row 2 - Patente / Motor / Chasis
$ws.Range("Y2").Value = "RGA001"
$ws.Range("Z2").Value = "ABC12SRGA001"
$ws.Range("AA2").Value = "ZAZ123SRGA001"

# Row 3 - Patente / Motor / Chasis
$ws.Range("Y3").Value = "RGA002"
$ws.Range("Z3").Value = "ABC12SRGA002"
$ws.Range("AA3").Value = "ZAZ123SRGA002"

# Row 4 - Patente / Motor / Chasis
$ws.Range("Y4").Value = "RGA003"
$ws.Range("Z4").Value = "ABC12SRGA003"
$ws.Range("AA4").Value = "ZAZ123SRGA003"

# SinAsistenciaMecanica for row 2
$ws.Range("AC2").Value = "Sí"

# Update selection to reflect where the user left off editing
$ws.Range("AC3").Select()

$wb.Save()
